$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11, columns A-G:
# A: employee_id, B: employee_name, C: department, D: absence_reason,
# E: absence_duration, F: absence_date, G: salary

$data = @(
    @(94251, "Carolina Lopes", "Vendas", "Consulta medica", 7, 45082, 9334.860000000001),
    @(17183, "Sr. Antony Sales", "P&D", "Viagem de negocios", 3, 45095, 8415.690000000001),
    @(57662, "Sr. João Pedro Silveira", "TI", "Problemas pessoais", 3, 45086, 5027.73),
    @(57159, "Maria Liz Melo", "Financeiro", "Viagem de negocios", 4, 45090, 9325.16),
    @(86497, "Emanuel Fernandes", "Recursos Humanos", "Viagem de negocios", 7, 45080, 3392.82),
    @(69918, "Isadora Silveira", "P&D", "Outros", 8, 45085, 5336.44),
    @(73558, "Srta. Gabrielly da Luz", "Recursos Humanos", "Problemas pessoais", 8, 45082, 3134.7),
    @(94144, "Maria Luísa Novais", "TI", "Consulta medica", 6, 45095, 9926.18),
    @(77062, "Juan Nogueira", "Recursos Humanos", "Consulta medica", 1, 45079, 5676.99),
    @(97686, "Sr. Thomas Siqueira", "Marketing", "Outros", 6, 45099, 9895.889999999999)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
